# Prise en compte des modifications lors de la conference video
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the training title
$ws.Range("A3").Value = "Formation : qh bon"

# Clear the company/trainer name next to "Trainer (s): "
$ws.Range("G3").Value = $null

# Replace the first attendee name, clear the rest of the attendee list
$ws.Range("A8").Value = "moi ahbo"
$ws.Range("A9:A18").Value = $null
